$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.182620048522949
$ws.Range("B1").Value = 1.752654075622559
$ws.Range("C1").Value = 6.944647789001465
$ws.Range("D1").Value = 1.624099731445312
$ws.Range("E1").Value = 0.9471976757049561
